$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "96586750-3/0"
$ws.Range("B2").Value = "CFINHRFLA"
$ws.Range("C2").Value = "L"
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 16956.9559
$ws.Range("F2").Value = "18/03/2025"
$ws.Range("G2").Value = "18/03/2025"
$ws.Range("H2").Value = 254354

# Update row 3
$ws.Range("A3").Value = "76513680-6/0"
$ws.Range("B3").Value = "CFINHRFLA"
$ws.Range("C3").Value = "L"
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 16956.9559
$ws.Range("F3").Value = "18/03/2025"
$ws.Range("G3").Value = "18/03/2025"
$ws.Range("H3").Value = 135656

# Remove row 4 entirely (was 96921130-0/0 row), shrinking used range to A1:H3
$ws.Range("A4:H4").ClearContents()
$ws.Rows.Item(4).Delete()
